$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''76.661.19'
$ws.Range('E2').Value = '  +0.14%  '
$ws.Range('D3').Value = '''2.931.99'
$ws.Range('E3').Value = '  +0.80%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = '''198.19'
$ws.Range('E5').Value = '  +0.29%  '
$ws.Range('D6').Value = '''595.40'
$ws.Range('E6').Value = '  -1.09%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('D8').Value = '''0.550'
$ws.Range('E8').Value = '  -1.33%  '
$ws.Range('D9').Value = '''0.199'
$ws.Range('E9').Value = '  +3.20%  '
$ws.Range('D10').Value = '''2.934.97'
$ws.Range('E10').Value = '  +1.00%  '
$ws.Range('D11').Value = '''0.441'
$ws.Range('E11').Value = '  +10.89%  '
$ws.Range('E12').Value = '  +0.31%  '
$ws.Range('D13').Value = '''3.476.16'
$ws.Range('E13').Value = '  +1.88%  '
$ws.Range('D14').Value = '''4.87'
$ws.Range('E14').Value = '  -1.76%  '
$ws.Range('D15').Value = '''28.37'
$ws.Range('E15').Value = '  +2.57%  '
$ws.Range('D16').Value = '''76.585.89'
$ws.Range('E16').Value = '  +0.18%  '
$ws.Range('D17').Value = '''0.0000190'
$ws.Range('E17').Value = '  -0.46%  '
$ws.Range('D18').Value = '''2.956.82'
$ws.Range('E18').Value = '  +2.06%  '
$ws.Range('D19').Value = '''13.49'
$ws.Range('E19').Value = '  +6.51%  '
$ws.Range('D20').Value = '''8.73'
$ws.Range('E20').Value = '  -3.53%  '
$ws.Range('D21').Value = '''373.93'
$ws.Range('E21').Value = '  -2.92%  '
$ws.Range('E22').Value = '  +3.21%  '
$ws.Range('E23').Value = '  -2.55%  '
$ws.Range('D24').Value = '''72.03'
$ws.Range('E24').Value = '  +0.20%  '
$ws.Range('E25').Value = '  -0.04%  '
$ws.Range('D26').Value = '''3.088.31'
$ws.Range('E26').Value = '  +1.83%  '
$ws.Range('D27').Value = '''4.26'
$ws.Range('E27').Value = '  -0.34%  '
$ws.Range('D28').Value = '''9.59'
$ws.Range('E28').Value = '  -2.37%  '
$ws.Range('E29').Value = '  -1.15%  '
$ws.Range('D30').Value = '''1.00'
$ws.Range('E30').Value = '  +0.28%  '
$ws.Range('D31').Value = '''8.32'
$ws.Range('E31').Value = '  +6.27%  '
$ws.Range('E32').Value = '  -3.90%  '
$ws.Range('D33').Value = '''501.10'
$ws.Range('E33').Value = '  -3.17%  '
$ws.Range('E34').Value = '  -0.09%  '
$ws.Range('E35').Value = '  +0.01%  '
$ws.Range('D36').Value = '''165.73'
$ws.Range('E36').Value = '  -0.61%  '
$ws.Range('B37').Value = 'EthereumClassic'
$ws.Range('C37').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D37').Value = '''20.13'
$ws.Range('E37').Value = '  -0.53%  '
$ws.Range('B38').Value = 'Cronos'
$ws.Range('C38').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D38').Value = '''0.111'
$ws.Range('E38').Value = '  +18.34%  '
$ws.Range('E39').Value = '  +11.88%  '
$ws.Range('D40').Value = '''19.96'
$ws.Range('E40').Value = '  +1.20%  '
$ws.Range('B41').Value = 'Kaspa'
$ws.Range('C41').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D41').Value = '''0.110'
$ws.Range('E41').Value = '  -6.07%  '
$ws.Range('B42').Value = 'USDe'
$ws.Range('C42').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D42').Value = '''1.00'
$ws.Range('E42').Value = '  +0.03%  '
$ws.Range('D43').Value = '''179.52'
$ws.Range('E43').Value = '  -2.36%  '
$ws.Range('D44').Value = '''4.91'
$ws.Range('E44').Value = '  -3.92%  '
$ws.Range('D45').Value = '''1.65'
$ws.Range('E45').Value = '  -2.89%  '
$ws.Range('D46').Value = '''40.05'
$ws.Range('E46').Value = '  -0.37%  '
$ws.Range('D48').Value = '''0.588'
$ws.Range('E48').Value = '  +0.36%  '
$ws.Range('E49').Value = '  +2.19%  '
$ws.Range('D50').Value = '''2.31'
$ws.Range('E50').Value = '  -3.46%  '
$ws.Range('D51').Value = '''22.36'
$ws.Range('E51').Value = '  +3.56%  '
